# "adding new fields in sale"
# Inserts a SEAL NO column ahead of VOL OBS and adds AMOUNT PAID / LOADING
# DATE / REMARKS columns after PAYMENT, extending the sale-log header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new column at I -- shifts old I:L (VOL OBS, VOL 20, SELLING
# PRICE, PAYMENT) one place to the right, into J:M, leaving I2 blank.
$ws.Range("I2").Insert(-4161)

# Give the new I2 header cell the same look as the other headers, then
# name it.
$ws.Range("J2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Value = "SEAL NO"

# Append the three brand-new trailing headers after PAYMENT (now M2),
# copying its header formatting first.
$ws.Range("M2").Copy()
$ws.Range("N2:P2").PasteSpecial(-4122)
$ws.Range("N2").Value = "AMOUNT PAID"
$ws.Range("O2").Value = "LOADING DATE"
$ws.Range("P2").Value = "REMARKS"

# Match the sheet's standard 20-wide columns for the newly introduced
# columns (I, N, O, P).
$refWidth = $ws.Columns(1).ColumnWidth
$ws.Columns(9).ColumnWidth = $refWidth
$ws.Columns(14).ColumnWidth = $refWidth
$ws.Columns(15).ColumnWidth = $refWidth
$ws.Columns(16).ColumnWidth = $refWidth

# Extend the title merge across the full new header width.
$ws.Range("A1:H1").UnMerge()
$ws.Range("A1:P1").Merge()

# Drop the obsolete truck-lookup validation that used to live on column D.
$ws.Range("D3:D100000").Validation.Delete()

# The YES/NO validation used to sit on the old SELLING PRICE column (which
# shifted to L); retarget it onto the PAYMENT column (M) instead.
$ws.Range("L3:L100000").Validation.Delete()
$ws.Range("M3:M100000").Validation.Add(3, 1, 1, '"YES, NO"')
$ws.Range("M3:M100000").Validation.IgnoreBlank = $false
